$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by duplicating "2021-Q4" (so it
#    inherits the exact same sheet/column formatting), inserting it
#    right before "2021-Q4". Final tab order becomes:
#      总计, 2022-Q3, 2021-Q4, 2020-Q4
# ---------------------------------------------------------------------
$q4_2021 = $wb.Worksheets.Item("2021-Q4")
$q4_2021.Copy($q4_2021)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q3"

# Overwrite the copied data with the new 2022-Q3 fund figures.
# A (index) and H (rank) are real numbers; B-G are kept as literal
# text so values such as "23.60" / "0.0010" keep their exact digits.
$newSheet.Range("A2").Value = 0
$txt2 = $newSheet.Range("B2:G2")
$txt2.NumberFormat = "@"
$newSheet.Range("B2").Value = "004402"
$newSheet.Range("C2").Value = "金信民旺债券C"
$newSheet.Range("D2").Value = "0.09"
$newSheet.Range("E2").Value = "23.60"
$newSheet.Range("F2").Value = "1.09"
$newSheet.Range("G2").Value = "0.0010"
$txt2.ClearFormats()
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$txt3 = $newSheet.Range("B3:G3")
$txt3.NumberFormat = "@"
$newSheet.Range("B3").Value = "004222"
$newSheet.Range("C3").Value = "金信民旺债券A"
$newSheet.Range("D3").Value = "0.08"
$newSheet.Range("E3").Value = "23.60"
$newSheet.Range("F3").Value = "1.09"
$newSheet.Range("G3").Value = "0.0009"
$txt3.ClearFormats()
$newSheet.Range("H3").Value = 10

# Restore the originally-active sheet selection (the Copy() above
# makes the new sheet the active tab, which we don't want).
$wb.Worksheets.Item("2020-Q4").Select()

# ---------------------------------------------------------------------
# 2) "总计" sheet: insert a new row 2 for "2022-Q3" and push the
#    existing "2021-Q4" / "2020-Q4" rows down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
